$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.182878228561681, 9.226618575922256, 3.082599426703578, 6.48142807727062, 0, 21.97352430845813)
    3 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
    4 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
    5 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 5.488907176552729)
    6 = @(0.1554434735375247, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1, 1.705647867635037)
    7 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 5.488907176552729)
    8 = @(0.06328177979961902, 1.65323645889881, 3.082599426703578, 6.48142807727062, 0, 11.28054574267263)
    9 = @(0.3464964993005633, 0.3375848360084654, 3.082599426703578, 6.48142807727062, 1, 10.24810883928323)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
